# arreglo cafee dea y bdea con trycatch
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows where column E ("cafee_bdea") had a stray "x" that must be cleared.
$eRowsToClear = @(2,3,4,5,6,7,8,9,10,14,15,16,38,39,40,50,51,52,53,54,55,56,57,59)
foreach ($r in $eRowsToClear) {
    $ws.Cells.Item($r, 5).Value = $null
}

# Rows where column D ("cafee_dea") was missing the "x" marker and must get it.
$dRowsToSet = 26..37
foreach ($r in $dRowsToSet) {
    $ws.Cells.Item($r, 4).Value = "x"
}

# Update the saved view state to match what was captured on commit.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 39
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G48").Select()
